$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text runs (in place, preserving shared string reuse) ---
# A8: "Volume 32   Number  9" -> "Volume 32   Number  10"  (replace trailing "9" at position 21)
$ws.Range("A8").Characters(21, 1).Text = "10"

# C9: "Report Covering the Week  2/24/2025  Through  3/2/2025"
#     -> "Report Covering the Week  3/3/2025  Through  3/9/2025"
# Replace the later occurrence first so the earlier offset stays valid.
$ws.Range("C9").Characters(47, 8).Text = "3/9/2025"
$ws.Range("C9").Characters(27, 9).Text = "3/3/2025"

# --- Update data grid (rows 15-31, columns C:N) ---

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("F15").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = 4

$ws.Range("D16").Value = 5

$ws.Range("E16").Value = -20

$ws.Range("F16").Value = 9

$ws.Range("G16").Value = 10

$ws.Range("H16").Value = -10

$ws.Range("I16").Value = 23

$ws.Range("J16").Value = 25

$ws.Range("K16").Value = -8

$ws.Range("L16").Value = 4.545454545454

$ws.Range("M16").Value = -11.538461538461

$ws.Range("N16").Value = -84.459459459459

$ws.Range("F15").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = 3

$ws.Range("D17").Value = 3

$ws.Range("E17").Value = 0

$ws.Range("F17").Value = 5

$ws.Range("G17").Value = 8

$ws.Range("H17").Value = -37.5

$ws.Range("I17").Value = 22

$ws.Range("J17").Value = 28

$ws.Range("K17").Value = -21.428571428571

$ws.Range("L17").Value = -29.032258064516

$ws.Range("M17").Value = 37.5

$ws.Range("N17").Value = -51.111111111111

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D18").PasteSpecial(-4122)

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E18").PasteSpecial(-4122)

$ws.Range("F18").Value = 5

$ws.Range("G18").Value = 4

$ws.Range("H18").Value = 25

$ws.Range("I18").Value = 13

$ws.Range("K18").Value = 8.333333333333

$ws.Range("L18").Value = -31.578947368421

$ws.Range("M18").Value = 8.333333333333

$ws.Range("N18").Value = -89.166666666666

$ws.Range("C19").Value = 7

$ws.Range("D19").Value = 8

$ws.Range("E19").Value = -12.5

$ws.Range("F19").Value = 19

$ws.Range("G19").Value = 27

$ws.Range("H19").Value = -29.629629629629

$ws.Range("I19").Value = 45

$ws.Range("J19").Value = 69

$ws.Range("K19").Value = -34.782608695652

$ws.Range("L19").Value = -39.189189189189

$ws.Range("M19").Value = -2.173913043478

$ws.Range("N19").Value = -62.184873949579

$ws.Range("F15").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = 4

$ws.Range("E20").Value = 300

$ws.Range("F20").Value = 4

$ws.Range("G20").Value = 4

$ws.Range("H20").Value = 0

$ws.Range("I20").Value = 5

$ws.Range("J20").Value = 6

$ws.Range("K20").Value = -16.666666666666

$ws.Range("L20").Value = -66.666666666666

$ws.Range("M20").Value = 150

$ws.Range("N20").Value = -94.444444444444

$ws.Range("C21").Value = 19

$ws.Range("D21").Value = 17

$ws.Range("E21").Value = 11.764705882352

$ws.Range("F21").Value = 43

$ws.Range("G21").Value = 54

$ws.Range("H21").Value = -20.37037037037

$ws.Range("I21").Value = 113

$ws.Range("J21").Value = 142

$ws.Range("K21").Value = -20.422535211267

$ws.Range("L21").Value = -30.674846625766

$ws.Range("M21").Value = 6.603773584905

$ws.Range("N21").Value = -78.679245283018

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("L22").Value = -42.857142857142

$ws.Range("C23").Value = 4

$ws.Range("D23").Value = 3

$ws.Range("E23").Value = 33.333333333333

$ws.Range("G23").Value = 12

$ws.Range("H23").Value = -33.333333333333

$ws.Range("I23").Value = 25

$ws.Range("J23").Value = 31

$ws.Range("K23").Value = -19.354838709677

$ws.Range("L23").Value = -13.793103448275

$ws.Range("M23").Value = 56.25

$ws.Range("C24").Value = 15

$ws.Range("D24").Value = 9

$ws.Range("E24").Value = 66.666666666666

$ws.Range("F24").Value = 48

$ws.Range("G24").Value = 35

$ws.Range("H24").Value = 37.142857142857

$ws.Range("I24").Value = 95

$ws.Range("J24").Value = 72

$ws.Range("K24").Value = 31.944444444444

$ws.Range("L24").Value = 14.457831325301

$ws.Range("M24").Value = 23.376623376623

$ws.Range("C25").Value = 1

$ws.Range("F15").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("D25").Value = 2

$ws.Range("M14").Copy()
$ws.Range("E25").PasteSpecial(-4122)
$ws.Range("E25").Value = -50

$ws.Range("F25").Value = 3

$ws.Range("H25").Value = -40

$ws.Range("I25").Value = 13

$ws.Range("J25").Value = 14

$ws.Range("K25").Value = -7.142857142857

$ws.Range("L25").Value = -58.064516129032

$ws.Range("C26").Value = 14

$ws.Range("D26").Value = 5

$ws.Range("E26").Value = 180

$ws.Range("F26").Value = 25

$ws.Range("G26").Value = 15

$ws.Range("H26").Value = 66.666666666666

$ws.Range("I26").Value = 53

$ws.Range("J26").Value = 52

$ws.Range("K26").Value = 1.923076923076

$ws.Range("L26").Value = 32.5

$ws.Range("M26").Value = 20.454545454545

$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("F15").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 1

$ws.Range("F15").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 1

$ws.Range("M14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = 0

$ws.Range("I28").Value = 5

$ws.Range("J28").Value = 3

$ws.Range("K28").Value = 66.666666666666

$ws.Range("L28").Value = -54.545454545454

$ws.Range("M14").Copy()
$ws.Range("M29").PasteSpecial(-4122)
$ws.Range("M29").Value = -100

$ws.Range("M14").Copy()
$ws.Range("M30").PasteSpecial(-4122)
$ws.Range("M30").Value = -100

$ws.Range("F15").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D31").Value = 2

$ws.Range("M14").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E31").Value = -100

$ws.Range("F15").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("G31").Value = 2

$ws.Range("M14").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("H31").Value = -100

$ws.Range("F15").Copy()
$ws.Range("J31").PasteSpecial(-4122)
$ws.Range("J31").Value = 2

$ws.Range("M14").Copy()
$ws.Range("K31").PasteSpecial(-4122)
$ws.Range("K31").Value = -50
